$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '68.105.23'
$ws.Range('E2').Value = '  +0.32%  '
$ws.Range('D3').Value = '3.885.88'
$ws.Range('E3').Value = '  -0.53%  '
$ws.Range('E4').Value = '  +0.32%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '483.35'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +1.11%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '145.48'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -1.93%  '
$ws.Range('E7').Value = '  +0.75%  '
$ws.Range('E8').Value = '  -0.09%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.739'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +2.06%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.176'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.69%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.0000355'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +0.80%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '42.91'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +0.95%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '10.58'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +2.98%  '
$ws.Range('D14').Value = '4.508.19'
$ws.Range('E14').Value = '  -0.24%  '
$ws.Range('D15').Value = '3.901.14'
$ws.Range('E15').Value = '  -1.91%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '14.30'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  -3.09%  '
$ws.Range('E17').Value = '  -0.56%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '20.21'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +1.62%  '
$ws.Range('E19').Value = '  +0.34%  '
$ws.Range('D20').Value = '68.116.49'
$ws.Range('E20').Value = '  +0.21%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '429.94'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -0.75%  '
$ws.Range('E22').Value = '  +6.04%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '14.85'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +3.01%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '89.53'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +2.48%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '11.78'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  +12.45%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '3.70'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +4.38%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '11.04'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +8.84%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '37.43'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  -1.85%  '
$ws.Range('E29').Value = '  -3.37%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '717.28'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -1.42%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '13.61'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +2.84%  '
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '2.92'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.72%  '
$ws.Range('D34').Value = '0.0₃0883'
$ws.Range('E34').Value = '  -2.49%  '
$ws.Range('B35').Value = 'NEARProtocol'
$ws.Range('C35').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '6.08'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  +12.66%  '
$ws.Range('B36').Value = 'InjectiveProtocol'
$ws.Range('C36').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '41.34'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  -1.51%  '
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '60.59'
$ws.Range('D37').Style = 'Normal'
$ws.Range('E37').Value = '  +3.00%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.397'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +15.74%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.145'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -4.00%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '0.995'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  -0.37%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '3.01'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +5.94%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '0.0495'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +5.67%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '3.10'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +4.07%  '
$ws.Range('E44').Value = '  -0.43%  '
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '0.144'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +2.46%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '3.35'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  +5.75%  '
$ws.Range('E47').Value = '  +0.21%  '
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '3.39'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -1.77%  '
$ws.Range('E49').Value = '  -2.54%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '145.21'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -0.30%  '
$ws.Range('E51').Value = '  -1.66%  '

Write-Host "Applied 92 cell updates"